$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "currency"
$ws.Range("D2").Value = "Dollar"
$ws.Range("D3").Value = "Pound"
$ws.Range("D4").Value = "Rupee"
$ws.Range("D5").Value = "Dollar"

$ws.Activate()
$ws.Range("B12").Select()
